$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.229372666666666
$ws.Range("H2").Value = 6.688117999999999
$ws.Range("I2").Value = 0.5889623983027473
$ws.Range("J2").Value = 0.5889623983027473
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.398836333333334
$ws.Range("N2").Value = 7.196509000000001
$ws.Range("O2").Value = 0.09386760623633866
$ws.Range("P2").Value = 0.09386760623633865
$ws.Range("Q2").Value = 5.347900153340222
$ws.Range("R2").Value = 48.131101380062
$ws.Range("S2").Value = 0.05528449049189194
$ws.Range("T2").Value = 0.05528449049189194

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.229372666666666
$ws.Range("H3").Value = 6.688117999999999
$ws.Range("I3").Value = 0.5889623983027473
$ws.Range("J3").Value = 0.5889623983027473
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.682092999999999
$ws.Range("N3").Value = 17.046279
$ws.Range("O3").Value = 0.2223430006085962
$ws.Range("P3").Value = 0.2223430006085962
$ws.Range("Q3").Value = 12.667502823658
$ws.Range("R3").Value = 114.007525412922
$ws.Range("S3").Value = 0.130951666884268
$ws.Range("T3").Value = 0.130951666884268

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.229372666666666
$ws.Range("H4").Value = 6.688117999999999
$ws.Range("I4").Value = 0.5889623983027473
$ws.Range("J4").Value = 0.5889623983027473
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.17451333333334
$ws.Range("N4").Value = 36.52354
$ws.Range("O4").Value = 0.4763944950360188
$ws.Range("P4").Value = 0.4763944950360188
$ws.Range("Q4").Value = 27.14152725530223
$ws.Range("R4").Value = 244.27374529772
$ws.Range("S4").Value = 0.2805784443346399
$ws.Range("T4").Value = 0.2805784443346399

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.229372666666666
$ws.Range("H5").Value = 6.688117999999999
$ws.Range("I5").Value = 0.5889623983027473
$ws.Range("J5").Value = 0.5889623983027473
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.300086333333334
$ws.Range("N5").Value = 15.900259
$ws.Range("O5").Value = 0.2073948981190463
$ws.Range("P5").Value = 0.2073948981190463
$ws.Range("Q5").Value = 11.81586760250689
$ws.Range("R5").Value = 106.342808422562
$ws.Range("S5").Value = 0.1221477965919475
$ws.Range("T5").Value = 0.1221477965919475

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.874264
$ws.Range("H6").Value = 2.622792
$ws.Range("I6").Value = 0.2309657016471988
$ws.Range("J6").Value = 0.2309657016471987
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.398836333333334
$ws.Range("N6").Value = 7.196509000000001
$ws.Range("O6").Value = 0.09386760623633866
$ws.Range("P6").Value = 0.09386760623633865
$ws.Range("Q6").Value = 2.097216248125334
$ws.Range("R6").Value = 18.874946233128
$ws.Range("S6").Value = 0.02168019753631893
$ws.Range("T6").Value = 0.02168019753631892

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.874264
$ws.Range("H7").Value = 2.622792
$ws.Range("I7").Value = 0.2309657016471988
$ws.Range("J7").Value = 0.2309657016471987
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.682092999999999
$ws.Range("N7").Value = 17.046279
$ws.Range("O7").Value = 0.2223430006085962
$ws.Range("P7").Value = 0.2223430006085962
$ws.Range("Q7").Value = 4.967649354551999
$ws.Range("R7").Value = 44.70884419096799
$ws.Range("S7").Value = 0.05135360714190797
$ws.Range("T7").Value = 0.05135360714190796

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.874264
$ws.Range("H8").Value = 2.622792
$ws.Range("I8").Value = 0.2309657016471988
$ws.Range("J8").Value = 0.2309657016471987
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 12.17451333333334
$ws.Range("N8").Value = 36.52354
$ws.Range("O8").Value = 0.4763944950360188
$ws.Range("P8").Value = 0.4763944950360188
$ws.Range("Q8").Value = 10.64373872485334
$ws.Range("R8").Value = 95.79364852368001
$ws.Range("S8").Value = 0.110030788806857
$ws.Range("T8").Value = 0.110030788806857

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.874264
$ws.Range("H9").Value = 2.622792
$ws.Range("I9").Value = 0.2309657016471988
$ws.Range("J9").Value = 0.2309657016471987
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.300086333333334
$ws.Range("N9").Value = 15.900259
$ws.Range("O9").Value = 0.2073948981190463
$ws.Range("P9").Value = 0.2073948981190463
$ws.Range("Q9").Value = 4.633674678125334
$ws.Range("R9").Value = 41.703072103128
$ws.Range("S9").Value = 0.04790110816211483
$ws.Range("T9").Value = 0.04790110816211482

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.136774
$ws.Range("H10").Value = 0.410322
$ws.Range("I10").Value = 0.0361333680411111
$ws.Range("J10").Value = 0.0361333680411111
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.398836333333334
$ws.Range("N10").Value = 7.196509000000001
$ws.Range("O10").Value = 0.09386760623633866
$ws.Range("P10").Value = 0.09386760623633865
$ws.Range("Q10").Value = 0.3280984406553333
$ws.Range("R10").Value = 2.952885965898
$ws.Range("S10").Value = 0.00339175276327572
$ws.Range("T10").Value = 0.00339175276327572

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.136774
$ws.Range("H11").Value = 0.410322
$ws.Range("I11").Value = 0.0361333680411111
$ws.Range("J11").Value = 0.0361333680411111
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.682092999999999
$ws.Range("N11").Value = 17.046279
$ws.Range("O11").Value = 0.2223430006085962
$ws.Range("P11").Value = 0.2223430006085962
$ws.Range("Q11").Value = 0.7771625879819998
$ws.Range("R11").Value = 6.994463291837999
$ws.Range("S11").Value = 0.008034001472355397
$ws.Range("T11").Value = 0.008034001472355397

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.136774
$ws.Range("H12").Value = 0.410322
$ws.Range("I12").Value = 0.0361333680411111
$ws.Range("J12").Value = 0.0361333680411111
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.17451333333334
$ws.Range("N12").Value = 36.52354
$ws.Range("O12").Value = 0.4763944950360188
$ws.Range("P12").Value = 0.4763944950360188
$ws.Range("Q12").Value = 1.665156886653333
$ws.Range("R12").Value = 14.98641197988
$ws.Range("S12").Value = 0.01721373762189574
$ws.Range("T12").Value = 0.01721373762189574

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.136774
$ws.Range("H13").Value = 0.410322
$ws.Range("I13").Value = 0.0361333680411111
$ws.Range("J13").Value = 0.0361333680411111
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.300086333333334
$ws.Range("N13").Value = 15.900259
$ws.Range("O13").Value = 0.2073948981190463
$ws.Range("P13").Value = 0.2073948981190463
$ws.Range("Q13").Value = 0.7249140081553332
$ws.Range("R13").Value = 6.524226073397999
$ws.Range("S13").Value = 0.007493876183584241
$ws.Range("T13").Value = 0.00749387618358424

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.544844
$ws.Range("H14").Value = 1.634532
$ws.Range("I14").Value = 0.1439385320089428
$ws.Range("J14").Value = 0.1439385320089428
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.398836333333334
$ws.Range("N14").Value = 7.196509000000001
$ws.Range("O14").Value = 0.09386760623633866
$ws.Range("P14").Value = 0.09386760623633865
$ws.Range("Q14").Value = 1.306991583198667
$ws.Range("R14").Value = 11.762924248788
$ws.Range("S14").Value = 0.01351116544485207
$ws.Range("T14").Value = 0.01351116544485207

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.544844
$ws.Range("H15").Value = 1.634532
$ws.Range("I15").Value = 0.1439385320089428
$ws.Range("J15").Value = 0.1439385320089428
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.682092999999999
$ws.Range("N15").Value = 17.046279
$ws.Range("O15").Value = 0.2223430006085962
$ws.Range("P15").Value = 0.2223430006085962
$ws.Range("Q15").Value = 3.095854278491999
$ws.Range("R15").Value = 27.862688506428
$ws.Range("S15").Value = 0.03200372511006482
$ws.Range("T15").Value = 0.03200372511006482

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.544844
$ws.Range("H16").Value = 1.634532
$ws.Range("I16").Value = 0.1439385320089428
$ws.Range("J16").Value = 0.1439385320089428
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.17451333333334
$ws.Range("N16").Value = 36.52354
$ws.Range("O16").Value = 0.4763944950360188
$ws.Range("P16").Value = 0.4763944950360188
$ws.Range("Q16").Value = 6.633210542586668
$ws.Range("R16").Value = 59.69889488328001
$ws.Range("S16").Value = 0.06857152427262612
$ws.Range("T16").Value = 0.06857152427262611

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.544844
$ws.Range("H17").Value = 1.634532
$ws.Range("I17").Value = 0.1439385320089428
$ws.Range("J17").Value = 0.1439385320089428
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.300086333333334
$ws.Range("N17").Value = 15.900259
$ws.Range("O17").Value = 0.2073948981190463
$ws.Range("P17").Value = 0.2073948981190463
$ws.Range("Q17").Value = 2.887720238198667
$ws.Range("R17").Value = 25.989482143788
$ws.Range("S17").Value = 0.02985211718139977
$ws.Range("T17").Value = 0.02985211718139977
